$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.408.39"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.891.72"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.690"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0736"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "2.166.93"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.723"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.909.42"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "35.387.67"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "0.0₃0822"
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").Value = "4.128.46"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0578"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -19.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0676"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "1.284.59"
$ws.Range("E44").Value = "  -3.93%  "
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0803"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.56%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.96%  "
